# Update the "Contacts" sheet: change the Title of the 3rd contact
# (row 4) from "Mrs." to "test", and update the selected cell to A5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contacts")
$ws.Activate()

$ws.Range("A4").Value = "test"

$ws.Range("A5").Select()
